$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Normalize font on existing data rows (3-38) - this is what drives the
#    font/style restructuring seen in the diff (new explicit Calibri font,
#    as opposed to the inherited theme font).
$ws.Range("A3:G38").Font.Name = "Calibri"

# 2. Add the new ticket row (row 39)
$ws.Range("A39").Value = "abertura-chamado"
$ws.Range("B39").Value = "'2025-03-19"
$ws.Range("C39").Value = "Maria Clara Rocha"
$ws.Range("E39").Value = "Pendente"
$ws.Range("F39").Value = "Lucas Rocha"
$ws.Range("G39").Value = "Média"
